# Auto-generated edit script applying numeric updates described in the commit diff.
$wb = $excel.ActiveWorkbook

# ALC row 2
$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H2").Value = 125.28571
$ws.Range("I2").Value = 125.28571
$ws.Range("K2").Value = 125.28571
$ws.Range("M2").Value = -12.28570999999999

# ALC row 17
$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H17").Value = 2500
$ws.Range("J17").Value = 2500
$ws.Range("L17").Value = 7500
$ws.Range("N17").Value = -7836

# ALC row 45
$ws = $wb.Worksheets.Item("ALC")
$ws.Range("J45").Value = 1017
$ws.Range("L45").Value = 3051
$ws.Range("N45").Value = -3435

# ALC row 132
$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H132").Value = 8488.143
$ws.Range("J132").Value = 6
$ws.Range("L132").Value = 18
$ws.Range("N132").Value = -5078

# ALC row 135
$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H135").Value = 2303
$ws.Range("I135").Value = 2360.5715
$ws.Range("K135").Value = 21245.1435
$ws.Range("M135").Value = -18710.1435

# ARM row 102
$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H102").Value = 8234.333000000001
$ws.Range("I102").Value = 8234.333000000001
$ws.Range("K102").Value = 8234.333000000001
$ws.Range("M102").Value = -6612.333000000001

# ARM row 122
$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H122").Value = 2701.8333
$ws.Range("I122").Value = 2701.8333
$ws.Range("K122").Value = 8105.499899999999
$ws.Range("M122").Value = -5655.499899999999

# BSM row 86
$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H86").Value = 5091.5713
$ws.Range("I86").Value = 5177.0835
$ws.Range("J86").Value = 4578.5
$ws.Range("K86").Value = 5177.0835
$ws.Range("L86").Value = 4578.5
$ws.Range("M86").Value = -4054.0835
$ws.Range("N86").Value = -6824.5

# BSM row 89
$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H89").Value = 5091.5713
$ws.Range("I89").Value = 5177.0835
$ws.Range("J89").Value = 4578.5
$ws.Range("K89").Value = 25885.4175
$ws.Range("L89").Value = 22892.5
$ws.Range("M89").Value = -20269.4175
$ws.Range("N89").Value = -34124.5

# BSM row 99
$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H99").Value = 947.6667
$ws.Range("I99").Value = 947.6667
$ws.Range("K99").Value = 947.6667
$ws.Range("M99").Value = 550.3333

# BSM row 105
$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H105").Value = 100996.5
$ws.Range("I105").Value = 1995
$ws.Range("K105").Value = 1995
$ws.Range("M105").Value = -248

# CRP row 45
$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H45").Value = 3800
$ws.Range("I45").Value = 0
$ws.Range("J45").Value = 3800
$ws.Range("K45").Value = 0
$ws.Range("L45").Value = 3800
$ws.Range("M45").ClearContents()
$ws.Range("N45").Value = -4986

# CRP row 58
$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H58").Value = 1855.5555
$ws.Range("I58").Value = 967
$ws.Range("J58").Value = 3632.6667
$ws.Range("K58").Value = 967
$ws.Range("L58").Value = 3632.6667
$ws.Range("M58").Value = -764
$ws.Range("N58").Value = -4038.6667

# CRP row 136
$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H136").Value = 1855.5555
$ws.Range("I136").Value = 967
$ws.Range("J136").Value = 3632.6667
$ws.Range("K136").Value = 2901
$ws.Range("L136").Value = 10898.0001
$ws.Range("M136").Value = -351
$ws.Range("N136").Value = -15998.0001

# CUL row 4
$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H4").Value = 1252.3334
$ws.Range("I4").Value = 556.1429000000001
$ws.Range("K4").Value = 1668.4287
$ws.Range("M4").Value = -1556.4287

# CUL row 69
$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H69").Value = 2604.6667
$ws.Range("J69").Value = 2604.6667
$ws.Range("L69").Value = 7814.000100000001
$ws.Range("N69").Value = -9436.000100000001

# CUL row 72
$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H72").Value = 2604.6667
$ws.Range("J72").Value = 2604.6667
$ws.Range("L72").Value = 23442.0003
$ws.Range("N72").Value = -31554.0003

# CUL row 88
$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H88").Value = 0
$ws.Range("J88").Value = 0
$ws.Range("L88").Value = 0
$ws.Range("N88").ClearContents()

# CUL row 91
$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H91").Value = 0
$ws.Range("J91").Value = 0
$ws.Range("L91").Value = 0
$ws.Range("N91").ClearContents()

# CUL row 131
$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H131").Value = 2626.6667
$ws.Range("I131").Value = 1700
$ws.Range("J131").Value = 2769.2307
$ws.Range("K131").Value = 5100
$ws.Range("L131").Value = 8307.6921
$ws.Range("M131").Value = -60
$ws.Range("N131").Value = -18387.6921

# GSM row 63
$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H63").Value = 0
$ws.Range("J63").Value = 0
$ws.Range("L63").Value = 0
$ws.Range("N63").ClearContents()

# GSM row 66
$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H66").Value = 0
$ws.Range("J66").Value = 0
$ws.Range("L66").Value = 0
$ws.Range("N66").ClearContents()

# GSM row 80
$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H80").Value = 15350
$ws.Range("I80").Value = 3125
$ws.Range("J80").Value = 21462.5
$ws.Range("K80").Value = 3125
$ws.Range("L80").Value = 21462.5
$ws.Range("M80").Value = -2127
$ws.Range("N80").Value = -23458.5

# GSM row 83
$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H83").Value = 15350
$ws.Range("I83").Value = 3125
$ws.Range("J83").Value = 21462.5
$ws.Range("K83").Value = 15625
$ws.Range("L83").Value = 107312.5
$ws.Range("M83").Value = -10633
$ws.Range("N83").Value = -117296.5

# LTW row 22
$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H22").Value = 774.75
$ws.Range("J22").Value = 949.75
$ws.Range("L22").Value = 949.75
$ws.Range("N22").Value = -1539.75

# LTW row 27
$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H27").Value = 774.75
$ws.Range("J27").Value = 949.75
$ws.Range("L27").Value = 949.75
$ws.Range("N27").Value = -1163.75

# LTW row 42
$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H42").Value = 40000000
$ws.Range("I42").Value = 0
$ws.Range("K42").Value = 0
$ws.Range("M42").ClearContents()

# LTW row 49
$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H49").Value = 40000000
$ws.Range("I49").Value = 0
$ws.Range("K49").Value = 0
$ws.Range("M49").ClearContents()

# LTW row 61
$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H61").Value = 0
$ws.Range("J61").Value = 0
$ws.Range("L61").Value = 0
$ws.Range("N61").ClearContents()

# LTW row 113
$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H113").Value = 0
$ws.Range("J113").Value = 0
$ws.Range("L113").Value = 0
$ws.Range("N113").ClearContents()

# LTW row 136
$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H136").Value = 5356.6665
$ws.Range("J136").Value = 0
$ws.Range("L136").Value = 0
$ws.Range("N136").ClearContents()

# WVR row 54
$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H54").Value = 0
$ws.Range("I54").Value = 0
$ws.Range("K54").Value = 0
$ws.Range("M54").ClearContents()

# WVR row 62
$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H62").Value = 10000
$ws.Range("J62").Value = 10000
$ws.Range("L62").Value = 10000
$ws.Range("N62").Value = -11248

# WVR row 65
$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H65").Value = 10000
$ws.Range("J65").Value = 10000
$ws.Range("L65").Value = 50000
$ws.Range("N65").Value = -56240

# WVR row 96
$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H96").Value = 995
$ws.Range("I96").Value = 995
$ws.Range("K96").Value = 995
$ws.Range("M96").Value = 378

# WVR row 100
$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H100").Value = 250
$ws.Range("I100").Value = 250
$ws.Range("K100").Value = 500
$ws.Range("M100").Value = 41

# WVR row 103
$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H103").Value = 20601
$ws.Range("J103").Value = 20601
$ws.Range("L103").Value = 20601
$ws.Range("N103").Value = -22945

# WVR row 122
$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H122").Value = 1499.6666
$ws.Range("I122").Value = 1499.6666
$ws.Range("K122").Value = 4498.9998
$ws.Range("M122").Value = -2048.9998
